$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds for row 8 (Ecuador - Liga Pro: Imbabura vs Tecnico U.) ---
$ws.Range("G8").Value = 2.55
$ws.Range("I8").Value = 2.8
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 3.6
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 2.25
$ws.Range("R8").Value = 1.62
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("W8").Value = 7.5
$ws.Range("Y8").Value = 10
$ws.Range("AC8").Value = 7.5
$ws.Range("AG8").Value = 351
$ws.Range("AL8").Value = 26
$ws.Range("AM8").Value = 41
$ws.Range("AV8").Value = 67
$ws.Range("AX8").Value = 17
$ws.Range("AY8").Value = 29
$ws.Range("BB8").Value = 251

# --- Update odds for row 9 (England - League One: Blackpool vs Wigan) ---
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 3.4
$ws.Range("K9").Value = 2.1
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 1.85
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 6
$ws.Range("AI9").Value = 17
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 12
$ws.Range("AU9").Value = 8
$ws.Range("BB9").Value = 201

# --- Remove the cancelled/dropped fixtures ---
# Row 10: ISRAEL - LIGAT HA'AL (Maccabi Haifa vs Hapoel Hadera) dropped entirely;
# every following row shifts up by one.
$ws.Rows(10).Delete()

# After the above deletion, the former row 15 (ROMANIA - LIGA 1, Farul Constanta vs
# Univ. Craiova) is now at row 14; drop it too so everything below shifts up again.
$ws.Rows(14).Delete()
